$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1608.6666
$ws.Range("I33").Value = 240
$ws.Range("J33").Value = 3524.8
$ws.Range("K33").Value = 240
$ws.Range("L33").Value = 3524.8
$ws.Range("M33").Value = -11
$ws.Range("N33").Value = -3982.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4240
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 5980
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 5980
$ws.Range("M116").Value = 942
$ws.Range("N116").Value = -12864

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2067.0588
$ws.Range("I125").Value = 2672
$ws.Range("K125").Value = 24048
$ws.Range("M125").Value = -21588

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2407.1226
$ws.Range("I127").Value = 369.33334
$ws.Range("J127").Value = 2540.0217
$ws.Range("K127").Value = 1108.00002
$ws.Range("L127").Value = 7620.0651
$ws.Range("M127").Value = 3851.99998
$ws.Range("N127").Value = -17540.0651

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 830.5417
$ws.Range("I135").Value = 767.1111
$ws.Range("J135").Value = 1020.8333
$ws.Range("K135").Value = 6903.9999
$ws.Range("L135").Value = 9187.4997
$ws.Range("M135").Value = -4368.9999
$ws.Range("N135").Value = -14257.4997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1011.2
$ws.Range("I137").Value = 1009.7222
$ws.Range("J137").Value = 1024.5
$ws.Range("K137").Value = 3029.1666
$ws.Range("L137").Value = 3073.5
$ws.Range("M137").Value = -479.1666
$ws.Range("N137").Value = -8173.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4103.7964
$ws.Range("I138").Value = 2029.6875
$ws.Range("J138").Value = 4875.558
$ws.Range("K138").Value = 6089.0625
$ws.Range("L138").Value = 14626.674
$ws.Range("M138").Value = -949.0625
$ws.Range("N138").Value = -24906.674

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2980.1035
$ws.Range("I141").Value = 2871
$ws.Range("K141").Value = 8613
$ws.Range("M141").Value = -3433

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2085.5715
$ws.Range("I61").Value = 1599.8334
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1599.8334
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1387.8334
$ws.Range("N61").Value = -5424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2185.2144
$ws.Range("I74").Value = 1285.2858
$ws.Range("J74").Value = 4885
$ws.Range("K74").Value = 1285.2858
$ws.Range("L74").Value = 4885
$ws.Range("M74").Value = -411.2858000000001
$ws.Range("N74").Value = -6633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2185.2144
$ws.Range("I77").Value = 1285.2858
$ws.Range("J77").Value = 4885
$ws.Range("K77").Value = 6426.429
$ws.Range("L77").Value = 24425
$ws.Range("M77").Value = -2058.429
$ws.Range("N77").Value = -33161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 49993
$ws.Range("J121").Value = 49993
$ws.Range("L121").Value = 49993
$ws.Range("N121").Value = -53487

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2620.4
$ws.Range("I122").Value = 2911.4
$ws.Range("J122").Value = 2038.4
$ws.Range("K122").Value = 8734.200000000001
$ws.Range("L122").Value = 6115.200000000001
$ws.Range("M122").Value = -6284.200000000001
$ws.Range("N122").Value = -11015.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1741.3214
$ws.Range("I132").Value = 1381.0869
$ws.Range("K132").Value = 4143.2607
$ws.Range("M132").Value = -1613.2607

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2085.5715
$ws.Range("I136").Value = 1599.8334
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4799.5002
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2249.5002
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20935.182
$ws.Range("I82").Value = 4169.5
$ws.Range("J82").Value = 30515.572
$ws.Range("K82").Value = 4169.5
$ws.Range("L82").Value = 30515.572
$ws.Range("M82").Value = -3786.5
$ws.Range("N82").Value = -31281.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 20935.182
$ws.Range("I85").Value = 4169.5
$ws.Range("J85").Value = 30515.572
$ws.Range("K85").Value = 4169.5
$ws.Range("L85").Value = 30515.572
$ws.Range("M85").Value = -2843.5
$ws.Range("N85").Value = -33167.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 139475.75
$ws.Range("I86").Value = 184499.83
$ws.Range("J86").Value = 4403.5
$ws.Range("K86").Value = 184499.83
$ws.Range("L86").Value = 4403.5
$ws.Range("M86").Value = -183376.83
$ws.Range("N86").Value = -6649.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 139475.75
$ws.Range("I89").Value = 184499.83
$ws.Range("J89").Value = 4403.5
$ws.Range("K89").Value = 922499.1499999999
$ws.Range("L89").Value = 22017.5
$ws.Range("M89").Value = -916883.1499999999
$ws.Range("N89").Value = -33249.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 22038.8
$ws.Range("J99").Value = 34671.332
$ws.Range("L99").Value = 34671.332
$ws.Range("N99").Value = -37667.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 22038.8
$ws.Range("J126").Value = 34671.332
$ws.Range("L126").Value = 104013.996
$ws.Range("N126").Value = -108953.996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 814.04
$ws.Range("J131").Value = 827.17206
$ws.Range("L131").Value = 2481.51618
$ws.Range("N131").Value = -12561.51618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2593
$ws.Range("I80").Value = 2411.4285
$ws.Range("K80").Value = 2411.4285
$ws.Range("M80").Value = -1413.4285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2593
$ws.Range("I83").Value = 2411.4285
$ws.Range("K83").Value = 12057.1425
$ws.Range("M83").Value = -7065.1425

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2653.0454
$ws.Range("I126").Value = 3303.4443
$ws.Range("J126").Value = 2202.7693
$ws.Range("K126").Value = 9910.332900000001
$ws.Range("L126").Value = 6608.3079
$ws.Range("M126").Value = -7440.332900000001
$ws.Range("N126").Value = -11548.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2410.348
$ws.Range("J132").Value = 3157.125
$ws.Range("L132").Value = 9471.375
$ws.Range("N132").Value = -14531.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2091.3635
$ws.Range("I7").Value = 1670.625
$ws.Range("K7").Value = 1670.625
$ws.Range("M7").Value = -1558.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1446580
$ws.Range("I46").Value = 683.3333
$ws.Range("J46").Value = 2531002.5
$ws.Range("K46").Value = 683.3333
$ws.Range("L46").Value = 2531002.5
$ws.Range("M46").Value = -495.3333
$ws.Range("N46").Value = -2531378.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 26964
$ws.Range("J98").Value = 26964
$ws.Range("L98").Value = 26964
$ws.Range("N98").Value = -32954

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1601.4
$ws.Range("I122").Value = 1333.6666
$ws.Range("J122").Value = 2003
$ws.Range("K122").Value = 4000.9998
$ws.Range("L122").Value = 6009
$ws.Range("M122").Value = -1550.9998
$ws.Range("N122").Value = -10909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2091.3635
$ws.Range("I126").Value = 1670.625
$ws.Range("K126").Value = 5011.875
$ws.Range("M126").Value = -2541.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1797.2142
$ws.Range("I136").Value = 1751.3334
$ws.Range("J136").Value = 1879.8
$ws.Range("K136").Value = 5254.0002
$ws.Range("L136").Value = 5639.4
$ws.Range("M136").Value = -2704.0002
$ws.Range("N136").Value = -10739.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 25935.4
$ws.Range("J41").Value = 25935.4
$ws.Range("L41").Value = 25935.4
$ws.Range("N41").Value = -26715.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2385.7144
$ws.Range("I122").Value = 1444.4445
$ws.Range("J122").Value = 4080
$ws.Range("K122").Value = 4333.333500000001
$ws.Range("L122").Value = 12240
$ws.Range("M122").Value = -1883.333500000001
$ws.Range("N122").Value = -17140
